# Fruta / hortaliza, semanal
# Update weekly Alcachofa price-report rows (34-43) on the active sheet.
# Each row's cells are set to the new values exactly as specified by the
# source diff (a weekly "shift" of the underlying records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("D34").Value = 45100
$ws.Range("H34").Value = "Argentina(o)"
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 16000
$ws.Range("L34").Value = 17000
$ws.Range("M34").Value = 16500
$ws.Range("N34").Value = "$/caja 50 unidades"
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 330
$ws.Range("Q34").Value = 50

# Row 35
$ws.Range("D35").Value = 44420
$ws.Range("J35").Value = 120
$ws.Range("K35").Value = 13000
$ws.Range("L35").Value = 14000
$ws.Range("M35").Value = 13500
$ws.Range("P35").Value = 338

# Row 36
$ws.Range("D36").Value = 44503
$ws.Range("J36").Value = 160

# Row 37
$ws.Range("D37").Value = 44505
$ws.Range("H37").Value = "Madrigal"
$ws.Range("J37").Value = 120
$ws.Range("K37").Value = 11000
$ws.Range("L37").Value = 12000
$ws.Range("M37").Value = 11500
$ws.Range("N37").Value = "$/caja 40 unidades"
$ws.Range("O37").Value = "Provincia del Elquí"
$ws.Range("P37").Value = 288
$ws.Range("Q37").Value = 40

# Row 38
$ws.Range("D38").Value = 44875
$ws.Range("H38").Value = "Española"
$ws.Range("J38").Value = 60
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = 10000
$ws.Range("N38").Value = "$/caja 30 unidades"
$ws.Range("O38").Value = "Provincia de Limarí"
$ws.Range("P38").Value = 333
$ws.Range("Q38").Value = 30

# Row 39
$ws.Range("D39").Value = 44488
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 11000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = 11500
$ws.Range("O39").Value = "Provincia del Elquí"
$ws.Range("P39").Value = 288

# Row 40
$ws.Range("D40").Value = 44426
$ws.Range("K40").Value = 13000
$ws.Range("L40").Value = 14000
$ws.Range("M40").Value = 13500
$ws.Range("O40").Value = "Región del Maule"
$ws.Range("P40").Value = 338

# Row 41
$ws.Range("D41").Value = 44510

# Row 42
$ws.Range("D42").Value = 44515
$ws.Range("J42").Value = 120

# Row 43
$ws.Range("D43").Value = 44490
$ws.Range("H43").Value = "Madrigal"
$ws.Range("J43").Value = 100
$ws.Range("K43").Value = 11000
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = 11500
$ws.Range("N43").Value = "$/caja 40 unidades"
$ws.Range("O43").Value = "Provincia del Elquí"
$ws.Range("P43").Value = 288
$ws.Range("Q43").Value = 40
